# 9th Stab - Cosmetic Changes
# Insert two new "week" columns (C and D) before the existing data column,
# shifting the old data column from C to E. Populate the two new columns
# with the same placeholder values as column B, and relabel the header row
# with the two newest week labels (Jun_15, Jun_17) while the old headers
# shift right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at C (existing C -> E)
$ws.Columns("C:D").Insert()

# Match the width of the newly inserted columns to the neighbouring column
# (character-width 7.14 round-trips to the same stored OOXML width, 8.0,
# that column E already carries)
$ws.Columns("C").ColumnWidth = 7.14
$ws.Columns("D").ColumnWidth = 7.14

$lastRow = $ws.UsedRange.Rows.Count

# Header row: B1 becomes the newest week, C1 the next newest, D1 keeps the
# previous B1 value ("Jun_13"), E1 keeps the old C1 value ("Jun_10").
$ws.Range("D1").Value = $ws.Range("B1").Value2
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# Data rows: columns C and D get the same placeholder value as column B.
for ($r = 2; $r -le $lastRow; $r++) {
    $placeholder = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 3).Value = $placeholder
    $ws.Cells.Item($r, 4).Value = $placeholder
}
